$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 206, shifting existing rows 206-272 down to 207-273
$ws.Rows.Item(206).Insert()

# Populate the new row 206 with its data (same static fields as the rest of
# the table, plus the new record's own date / volume / price values)
$ws.Range("A206").Value = 3
$ws.Range("B206").Value = "Femacal de La Calera"
$ws.Range("C206").Value = "Coquimbo"
$ws.Range("D206").Value = 44588
$ws.Range("E206").Value = 5
$ws.Range("F206").Value = 100112012
$ws.Range("G206").Value = "Espinaca"
$ws.Range("H206").Value = "Sin especificar"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 160
$ws.Range("K206").Value = 3500
$ws.Range("L206").Value = 4000
$ws.Range("M206").Value = 3750
$ws.Range("N206").Value = "$/docena de atados (3 kilos)"
$ws.Range("O206").Value = "Provincia de Quillota"
$ws.Range("P206").Value = 1250
$ws.Range("Q206").Value = 3
$ws.Range("R206").Value = "Hortaliza"

# Match the date-column number format used by the rest of column D
$ws.Range("D206").NumberFormat = $ws.Range("D207").NumberFormat
